# Form the consolidated report: set the "Absent" (H) column values so
# that they correctly reflect the inverse of the "Real" (E) column for
# each attendance date row, filling in the cells that were left blank
# and fixing the ones that were miscalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
